$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at row 2, pushing existing data down
$ws.Range("A2:A3").EntireRow.Insert()

# Copy formatting from the (now shifted) row 4 - former row 2 - down onto the
# two newly inserted blank rows so they match the rest of the table's style.
$ws.Range("A4:J4").Copy()
$ws.Range("A2:J3").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# New row 2: ChatGPT
$ws.Range("A2").Value = "ChatGPT"
$ws.Range("B2").Value = 45888
$ws.Range("C2").Value = 0.88
$ws.Range("D2").Value = 0.79
$ws.Range("E2").Value = 0.771
$ws.Range("F2").Value = 0.8179999999999999
$ws.Range("G2").Value = 0.8139999999999999
$ws.Range("H2").Value = 1895
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = "model: ChatGPT ThreeMainClassModel"

# New row 3: Lexical search
$ws.Range("A3").Value = "Lexical search"
$ws.Range("B3").Value = 45510
$ws.Range("C3").Value = 0.8389999866485596
$ws.Range("D3").Value = 0.7300000190734863
$ws.Range("E3").Value = 0.7089999914169312
$ws.Range("F3").Value = 0.7699999809265137
$ws.Range("G3").Value = 0.5789999961853027
$ws.Range("H3").Value = 1895
$ws.Range("I3").Value = 7578
$ws.Range("J3").Value = "Ngram 1"
